$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.29"
$ws.Range("E2").Value = "'1.63%"
$ws.Range("D3").Value = "'36.21"
$ws.Range("E3").Value = "'-0.69%"
$ws.Range("D4").Value = "'5.076"
$ws.Range("E4").Value = "'1.88%"
$ws.Range("D5").Value = "'0.07934"
$ws.Range("E5").Value = "'2.58%"
$ws.Range("D6").Value = "'2.166"
$ws.Range("E6").Value = "'4.65%"
$ws.Range("D7").Value = "'8.025"
$ws.Range("E7").Value = "'1.54%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9315"
$ws.Range("E8").Value = "'1.02%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09892"
$ws.Range("E9").Value = "'1.62%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1873"
$ws.Range("E10").Value = "'1.00%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09017"
$ws.Range("E11").Value = "'5.30%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03627"
$ws.Range("E12").Value = "'3.11%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09934"
$ws.Range("E13").Value = "'-0.07%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001442"
$ws.Range("E14").Value = "'-1.80%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005710"
$ws.Range("E15").Value = "'1.46%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.448"
$ws.Range("E16").Value = "'-0.38%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.167"
$ws.Range("E17").Value = "'3.59%"
$ws.Range("E18").Value = "'15.44%"
$ws.Range("E19").Value = "'-1.03%"
$ws.Range("D20").Value = "'0.1357"
$ws.Range("E20").Value = "'1.22%"
$ws.Range("D21").Value = "'5.083"
$ws.Range("E21").Value = "'6.63%"
$ws.Range("D22").Value = "'0.2192"
$ws.Range("E22").Value = "'-0.17%"
$ws.Range("D23").Value = "'0.04583"
$ws.Range("E23").Value = "'-0.26%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'1.00%"
$ws.Range("D25").Value = "'0.004765"
$ws.Range("E25").Value = "'-6.26%"
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'-6.80%"
$ws.Range("D39").Value = "'0.01944"
$ws.Range("E39").Value = "'10.51%"
$ws.Range("D40").Value = "'0.04922"
$ws.Range("E40").Value = "'6.07%"
$ws.Range("D41").Value = "'0.007805"
$ws.Range("E41").Value = "'4.98%"
$ws.Range("D42").Value = "'0.1393"
$ws.Range("E42").Value = "'0.34%"
$ws.Range("D43").Value = "'0.007734"
$ws.Range("E43").Value = "'0.38%"
$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'-6.47%"
$ws.Range("D45").Value = "'0.01144"
$ws.Range("E45").Value = "'10.91%"
$ws.Range("D46").Value = "'0.00006227"
$ws.Range("E46").Value = "'0.95%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.39%"
$ws.Range("D48").Value = "'51.96"
$ws.Range("E48").Value = "'45.63%"
$ws.Range("D49").Value = "'0.001804"
$ws.Range("E49").Value = "'-9.63%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.39%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.39%"
